# Update moeda (currency) quotes and "last updated" timestamps.
# Source data changed from "11 de jun." readings to "16 de jun." readings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Dólar (USD)
$ws.Range("E2").Value = "5,49"
$ws.Range("F2").Value = "16 de jun., 19:31 UTC ·"

# Row 3 - Euro (EUR)
$ws.Range("E3").Value = "6,36"
$ws.Range("F3").Value = "16 de jun., 19:26 UTC ·"

# Row 4 - Real (BRL)
$ws.Range("F4").Value = "16 de jun., 19:31 UTC ·"

# Row 5 - Libra Esterlina (GBP)
$ws.Range("E5").Value = "7,47"
$ws.Range("F5").Value = "16 de jun., 19:30 UTC ·"

# Row 6 - Iene (JPY)
$ws.Range("F6").Value = "16 de jun., 19:30 UTC ·"

# Row 7 - Franco Suíço (CHF)
$ws.Range("F7").Value = "16 de jun., 19:26 UTC ·"

# Row 8 - Dólar Australiano (AUD)
$ws.Range("E8").Value = "3,59"
$ws.Range("F8").Value = "16 de jun., 19:31 UTC ·"

# Row 9 - Peso Mexicano (MXN)
$ws.Range("F9").Value = "16 de jun., 19:31 UTC ·"

# Row 10 - Dólar Canadiano (CAD)
$ws.Range("E10").Value = "4,05"
$ws.Range("F10").Value = "16 de jun., 19:31 UTC ·"

# Row 11 - Dólar de Hong Kong (HKD)
$ws.Range("E11").Value = "0,70"
$ws.Range("F11").Value = "16 de jun., 19:31 UTC ·"

# Row 12 - Yuan Chinês (CNY)
$ws.Range("F12").Value = "16 de jun., 19:29 UTC ·"

# Row 13 - Rúpia Indiana (INR)
# (Force text parsing: "0,064" has exactly 3 digits after the comma, which
# the COM layer otherwise misreads as a grouped integer "64".)
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0,064"
$ws.Range("E13").ClearFormats()
$ws.Range("F13").Value = "16 de jun., 19:30 UTC ·"

# Row 14 - Peso Chileno (CLP)
$ws.Range("F14").Value = "16 de jun., 11:19 UTC ·"

# Row 15 - Peso Argentino (ARS)
$ws.Range("E15").Value = "0,0046"
$ws.Range("F15").Value = "16 de jun., 19:30 UTC ·"

# Row 16 - Peso Colombiano (COP)
$ws.Range("F16").Value = "16 de jun., 19:31 UTC ·"

# Row 17 - Rúpia Russa (RUB)
$ws.Range("E17").Value = "1,09"
$ws.Range("F17").Value = "16 de jun., 19:30 UTC ·"

# Row 18 - Riyal Saudi (SAR)
$ws.Range("E18").Value = "1,46"
$ws.Range("F18").Value = "16 de jun., 19:30 UTC ·"

# Row 19 - Dólar de Singapura (SGD)
$ws.Range("E19").Value = "4,29"
$ws.Range("F19").Value = "16 de jun., 19:31 UTC ·"

# Row 20 - Peso Filipino (PHP)
# (Same 3-decimal-digit text-vs-number quirk as row 13.)
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0,097"
$ws.Range("E20").ClearFormats()
$ws.Range("F20").Value = "16 de jun., 19:30 UTC ·"

# Row 21 - Yuan de Taiwan (TWD)
$ws.Range("E21").Value = "4,11"
$ws.Range("F21").Value = "16 de jun., 19:31 UTC ·"

# Row 22 - Dinar Iraquiano (IQD)
$ws.Range("F22").Value = "16 de jun., 19:30 UTC ·"

# Row 23 - Rúpia Sri Lanka (LKR)
$ws.Range("F23").Value = "16 de jun., 19:29 UTC ·"

# Row 24 - Yuan Chinês (CNY)
$ws.Range("F24").Value = "16 de jun., 19:29 UTC ·"

# Row 25 - Won Sul-Coreano (KRW)
$ws.Range("F25").Value = "16 de jun., 19:31 UTC ·"
